{"js": "// Replace the twenty-five division problems in the practice table with the\n// new set of problems, leaving everything else (fonts, sizes, layout, the\n// date line, empty answer rows, etc.) untouched.\nconst replacements = [\n  [\"237\u00f78=\", \"656\u00f73=\"],\n  [\"612\u00f75=\", \"208\u00f75=\"],\n  [\"972\u00f77=\", \"658\u00f73=\"],\n  [\"916\u00f74=\", \"191\u00f77=\"],\n  [\"555\u00f75=\", \"757\u00f72=\"],\n  [\"585\u00f74=\", \"626\u00f74=\"],\n  [\"984\u00f77=\", \"230\u00f78=\"],\n  [\"120\u00f72=\", \"242\u00f73=\"],\n  [\"455\u00f77=\", \"738\u00f73=\"],\n  [\"423\u00f72=\", \"523\u00f76=\"],\n  [\"939\u00f74=\", \"483\u00f75=\"],\n  [\"585\u00f76=\", \"152\u00f74=\"],\n  [\"734\u00f73=\", \"943\u00f76=\"],\n  [\"202\u00f77=\", \"892\u00f72=\"],\n  [\"490\u00f73=\", \"958\u00f75=\"],\n  [\"230\u00f75=\", \"871\u00f75=\"],\n  [\"313\u00f78=\", \"216\u00f73=\"],\n  [\"733\u00f74=\", \"659\u00f74=\"],\n  [\"929\u00f74=\", \"950\u00f74=\"],\n  [\"177\u00f75=\", \"485\u00f73=\"],\n  [\"259\u00f78=\", \"569\u00f72=\"],\n  [\"512\u00f76=\", \"786\u00f72=\"],\n  [\"772\u00f76=\", \"134\u00f75=\"],\n  [\"355\u00f72=\", \"151\u00f75=\"],\n  [\"356\u00f79=\", \"587\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the twenty-five division problems in the practice table with the\n# new set of problems, leaving everything else (fonts, sizes, layout, the\n# date line, empty answer rows, etc.) untouched.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"237\u00f78=\", \"656\u00f73=\"),\n    @(\"612\u00f75=\", \"208\u00f75=\"),\n    @(\"972\u00f77=\", \"658\u00f73=\"),\n    @(\"916\u00f74=\", \"191\u00f77=\"),\n    @(\"555\u00f75=\", \"757\u00f72=\"),\n    @(\"585\u00f74=\", \"626\u00f74=\"),\n    @(\"984\u00f77=\", \"230\u00f78=\"),\n    @(\"120\u00f72=\", \"242\u00f73=\"),\n    @(\"455\u00f77=\", \"738\u00f73=\"),\n    @(\"423\u00f72=\", \"523\u00f76=\"),\n    @(\"939\u00f74=\", \"483\u00f75=\"),\n    @(\"585\u00f76=\", \"152\u00f74=\"),\n    @(\"734\u00f73=\", \"943\u00f76=\"),\n    @(\"202\u00f77=\", \"892\u00f72=\"),\n    @(\"490\u00f73=\", \"958\u00f75=\"),\n    @(\"230\u00f75=\", \"871\u00f75=\"),\n    @(\"313\u00f78=\", \"216\u00f73=\"),\n    @(\"733\u00f74=\", \"659\u00f74=\"),\n    @(\"929\u00f74=\", \"950\u00f74=\"),\n    @(\"177\u00f75=\", \"485\u00f73=\"),\n    @(\"259\u00f78=\", \"569\u00f72=\"),\n    @(\"512\u00f76=\", \"786\u00f72=\"),\n    @(\"772\u00f76=\", \"134\u00f75=\"),\n    @(\"355\u00f72=\", \"151\u00f75=\"),\n    @(\"356\u00f79=\", \"587\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
